$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with changed quantity/stats values (rows 4-11, 15-16 per diff)
$ws.Range("G4").Value = -1349
$ws.Range("H4").Value = 1.06
$ws.Range("I4").Value = 0.28

$ws.Range("G5").Value = -1349
$ws.Range("H5").Value = 1.06
$ws.Range("I5").Value = 0.28

$ws.Range("G6").Value = -1349
$ws.Range("H6").Value = 1.06
$ws.Range("I6").Value = 0.28

# Row 7 and 8 swap their E/F values (id_produto/produto) plus other stats
$ws.Range("E7").Value = 13544
$ws.Range("F7").Value = "MOUSE SEM FIO 3 BOTOES 1000DPI COLOR FIT BRANCO 1709 R8"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 1.07
$ws.Range("I7").Value = 0.26

$ws.Range("E8").Value = 13244
$ws.Range("F8").Value = "MOUSE SEM FIO 3 BOTOES 1000DPI COLOR FIT AZUL 1709 R8"
$ws.Range("G8").Value = -15
$ws.Range("H8").Value = 1.06
$ws.Range("I8").Value = 0.25

$ws.Range("G9").Value = -99
$ws.Range("H9").Value = 1.04
$ws.Range("I9").Value = 0.19

$ws.Range("G10").Value = -56
$ws.Range("I10").Value = 0.13

$ws.Range("G11").Value = -1349
$ws.Range("H11").Value = 1.06
$ws.Range("I11").Value = 0.28

$ws.Range("G15").Value = -129
$ws.Range("I15").Value = 0.21

$ws.Range("G16").Value = -5
$ws.Range("I16").Value = 0.28

# Append new rows 18-25
$newRows = @(
    @("2025-08-19", 3, "BEMOL S/A", "401315", 13854, "CANETA APAGAVEL GEL 0.7 AZUL/PRETA - CORES SORTIDAS JOCAR OFFICE", 0, 1.24, 0.5600000000000001),
    @("2025-08-19", 2, "BEMOL S/A", "401319", 12945, "FONE DE OUVIDO SEM FIO BT BASIKE FON-9856", -99, 1.04, 0.19),
    @("2025-08-19", 2, "BEMOL S/A", "401332", 10525, "BOMBA AUTOMATICA PARA GALAO DE AGUA RECARREGAVEL USB", -129, 1.04, 0.21),
    @("2025-08-19", 2, "BEMOL S/A", "401336", 13977, "SAPATEIRA MATERIAL PP, FERRO E TNT CAPACIDADE18 PARES, SUPORTA ATE 15KG", -38, 1.03, 0.16),
    @("2025-08-19", 2, "BEMOL S/A", "401337", 10130, "FONE DE OUVIDO SEM FIO A GOLD V5.3", -1349, 1.06, 0.28),
    @("2025-08-19", 2, "BEMOL S/A", "401348", 10130, "FONE DE OUVIDO SEM FIO A GOLD V5.3", -1349, 1.06, 0.28),
    @("2025-08-19", 3, "BEMOL S/A", "401363", 10130, "FONE DE OUVIDO SEM FIO A GOLD V5.3", -1349, 1.06, 0.28),
    @("2025-08-19", 2, "BEMOL S/A", "401365", 10130, "FONE DE OUVIDO SEM FIO A GOLD V5.3", -1349, 1.06, 0.28)
)

$r = 18
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}
